$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '63.106.45'
$c.Style = "Normal"
$ws.Range('E2').Value = '  -0.91%  '

$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '3.151.35'
$c.Style = "Normal"
$ws.Range('E3').Value = '  +1.09%  '

$ws.Range('E4').Value = '  -0.10%  '

$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '588.17'
$c.Style = "Normal"
$ws.Range('E5').Value = '  -1.77%  '

$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '138.09'
$c.Style = "Normal"
$ws.Range('E6').Value = '  -3.06%  '

$ws.Range('E7').Value = '  -0.07%  '

$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '3.148.31'
$c.Style = "Normal"
$ws.Range('E8').Value = '  +1.43%  '

$ws.Range('E9').Value = '  -0.05%  '

$ws.Range('E10').Value = '  -0.94%  '

$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '5.29'
$c.Style = "Normal"
$ws.Range('E11').Value = '  -0.62%  '

$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '0.458'
$c.Style = "Normal"
$ws.Range('E12').Value = '  -1.21%  '

$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '0.0000245'
$c.Style = "Normal"
$ws.Range('E13').Value = '  -2.09%  '

$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '34.19'
$c.Style = "Normal"
$ws.Range('E14').Value = '  -2.16%  '

$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '3.670.79'
$c.Style = "Normal"

$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '0.120'
$c.Style = "Normal"
$ws.Range('E16').Value = '  +1.15%  '

$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '3.150.36'
$c.Style = "Normal"
$ws.Range('E17').Value = '  +0.99%  '

$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '63.079.73'
$c.Style = "Normal"
$ws.Range('E18').Value = '  -1.10%  '

$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '6.66'
$c.Style = "Normal"
$ws.Range('E19').Value = '  -1.05%  '

$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '476.36'
$c.Style = "Normal"
$ws.Range('E20').Value = '  -0.62%  '

$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '14.02'
$c.Style = "Normal"
$ws.Range('E21').Value = '  -3.36%  '

$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '0.702'
$c.Style = "Normal"
$ws.Range('E22').Value = '  +0.05%  '

$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '7.73'
$c.Style = "Normal"
$ws.Range('E23').Value = '  +2.27%  '

$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '84.60'
$c.Style = "Normal"
$ws.Range('E24').Value = '  -2.71%  '

$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '12.99'
$c.Style = "Normal"
$ws.Range('E25').Value = '  -1.71%  '

$ws.Range('E27').Value = '  -0.66%  '

$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '7.09'
$c.Style = "Normal"
$ws.Range('E28').Value = '  +0.90%  '

$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '7.94'
$c.Style = "Normal"
$ws.Range('E29').Value = '  -2.89%  '

$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '2.11'
$c.Style = "Normal"
$ws.Range('E30').Value = '  +3.85%  '

$ws.Range('E31').Value = '  -0.01%  '

$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '26.85'
$c.Style = "Normal"
$ws.Range('E32').Value = '  -0.53%  '

$ws.Range('E33').Value = '  -4.13%  '

$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '2.54'
$c.Style = "Normal"
$ws.Range('E34').Value = '  -3.21%  '

$ws.Range('E35').Value = '  -2.19%  '

$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '5.81'
$c.Style = "Normal"
$ws.Range('E36').Value = '  -2.53%  '

$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '52.44'
$c.Style = "Normal"
$ws.Range('E37').Value = '  -0.07%  '

$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '0.0₃0698'
$c.Style = "Normal"
$ws.Range('E38').Value = '  -6.89%  '

$ws.Range('E39').Value = '  -0.06%  '

$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '416.04'
$c.Style = "Normal"
$ws.Range('E40').Value = '  -4.39%  '

$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '2.76'
$c.Style = "Normal"
$ws.Range('E41').Value = '  -5.63%  '

$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '8.28'
$c.Style = "Normal"
$ws.Range('E42').Value = '  +0.74%  '

$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '2.927.77'
$c.Style = "Normal"
$ws.Range('E43').Value = '  +2.81%  '

$ws.Range('E44').Value = '  -6.11%  '

$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '0.262'
$c.Style = "Normal"
$ws.Range('E45').Value = '  +1.71%  '

$ws.Range('E46').Value = '  +0.04%  '

$ws.Range('E47').Value = '  -2.90%  '

$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '25.41'
$c.Style = "Normal"
$ws.Range('E48').Value = '  -0.63%  '

$ws.Range('E49').Value = '  +0.44%  '

$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '2.24'
$c.Style = "Normal"
$ws.Range('E50').Value = '  -7.50%  '

$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '120.94'
$c.Style = "Normal"
$ws.Range('E51').Value = '  -0.45%  '
